# Adicionado o paramentro para enviar a data e hora da requisicao
# Adds five new customer rows (23-27) to the "Clientes" sheet, mirroring the
# existing "BRUNO DE FRAGA" test rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A = "BRUNO DE FRAGA"; B = "123123"; C = "023.709.450-95"; E = "92320-195"; F = "bruno@gmail.com"; G = "51989043802"; H = "Rua 3 Pinheiros I, 27" },
    @{ A = "BRUNO DE FRAGA"; B = "123123"; C = "3123123";        E = "92320-195"; F = "bruno@gmail.com"; G = "51989043802"; H = "Rua 3 Pinheiros I, 27" },
    @{ A = "BRUNO DE FRAGA"; B = "123123"; C = "123123";         E = "92320-195"; F = "bruno@gmail.com"; G = "51989043802"; H = "Rua 3 Pinheiros I, 27" },
    @{ A = "BRUNO DE FRAGA"; B = "123123"; C = "023.709.450-95"; E = "92320-195"; F = "bruno@gmail.com"; G = "51989043802"; H = "Rua 3 Pinheiros I, 27" },
    @{ A = "BRUNO DE FRAGA"; B = "123123"; C = "4123123";        E = "92320-195"; F = "bruno@gmail.com"; G = "51989043802"; H = "Rua 3 Pinheiros I, 27" }
)

$startRow = 23
$cols = @("A", "B", "C", "E", "F", "G", "H")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    foreach ($col in $cols) {
        $cellRef = $col + $r
        $val = $row[$col]
        # Keep digit-only values (CPF/phone/password placeholders) as text
        # instead of letting them be auto-converted to numbers, just like
        # the source data that already lives in rows 2-22.
        if ($val -match '^[0-9]+$') {
            $ws.Range($cellRef).NumberFormat = "@"
        }
        $ws.Range($cellRef).Value = $val
    }
}

Write-Output "Added rows 23-27 to Clientes sheet"
